# Add a new "Adjusted2" worksheet (sheetId 4) after the existing "Adjusted"
# sheet, containing theta/Jxx/Jyy/beta/gamma/trace data for the
# "stokes parameter input" variant of the polarimeter calculation.

$wb = $excel.ActiveWorkbook
$srcWs = $wb.Worksheets.Item("Adjusted")

# Insert the new sheet right after the last existing sheet so it lands at
# the end of the tab strip (Worksheets.Add with no args inserts at the
# front, so we explicitly say "after the current last sheet").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Adjusted2"

# Copy the header row (value + formatting) from the existing "Adjusted"
# sheet so the new sheet's header reuses the exact same style (bold,
# thin border, centered/top-aligned) instead of minting new style records.
$srcWs.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$srcWs.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows: theta, Jxx, Jyy, beta, gamma, trace
$data = @(
    @(0, 0.0004618726657399049, 0.9995381342511355, -0.004801445775495194, 0, 0.9991228029116703),
    @(5, 0.0315759714091005, 0.9684240233729264, -0.1717474004983522, 0, 0.9978364701721175),
    @(10, 0.121417889823968, 0.8785820704600988, -0.3217453404712551, 0, 0.9936888867331871),
    @(15, 0.251741460114978, 0.7482584437466685, -0.4269638638656761, 0, 0.9878607434732226),
    @(20, 0.4241549422155633, 0.5758449075404576, -0.4855899895252066, 0, 0.9831000484003465),
    @(25, 0.5939256348482375, 0.4060741882444195, -0.4747920702346514, 0, 0.9684989260036585),
    @(30, 0.7437721743770547, 0.2562276876743538, -0.4160445543303309, 0, 0.9650358176843653),
    @(35, 0.8718430725792202, 0.1281568523263641, -0.3046078503604101, 0, 0.9621064070049569),
    @(40, 0.9506213267939327, 0.04937863873068832, -0.1554611713297683, 0, 0.9544555085008001),
    @(45, 0.9778494965890689, 0.02215048487196899, 0.01034466067816938, 0, 0.9568943059686518),
    @(50, 0.9454368594929669, 0.05456310445527547, 0.1743640845519795, 0, 0.957633655619021),
    @(55, 0.8668626798317627, 0.1331372439452539, 0.3131274581478746, 0, 0.9652740415027411),
    @(60, 0.7362078201682305, 0.2637920477012068, 0.4238019795912741, 0, 0.9708044347182189),
    @(65, 0.58080990058667, 0.4191899350227959, 0.4819649903581523, 0, 0.9776408461057811),
    @(70, 0.4154959398797902, 0.584503902238563, 0.4830881369867478, 0, 0.9810299839833514),
    @(75, 0.25485053654043, 0.7451493665712157, 0.428768891175602, 0, 0.9878818985563395),
    @(80, 0.1186921997024051, 0.8813077615318918, 0.318567373911987, 0, 0.9937615522489089),
    @(85, 0.03335612927895004, 0.966643863921934, 0.1758824771011767, 0, 0.9973822825208924),
    @(90, 0.0006876719185082963, 0.9993123347381293, -0.00240164221016018, 0, 0.9986371510230498)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Restore the originally active sheet so this edit doesn't change which
# tab is selected when the workbook is opened.
$wb.Worksheets.Item(1).Activate()

Write-Output "Adjusted2 sheet created"
